# Update the date placeholder on the first slide (title slide) from
# 2017-03-26 to 2017-04-28, per the commit's "Regenerated pptx" pass.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$dateShape = $s.Shapes.Item("Date 3")
$dateShape.TextFrame.TextRange.Text = "2017-04-28"
